$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ingredients list for the "Ivy" cocktail row: drop the "dry " qualifier on gin.
$ws.Range("B2").Value = "gin, green chartreuse, dry vermouth, absinthe, orange bitters"

# Move the active selection to B3, as recorded in the saved view state.
$ws.Range("B3").Select()
